$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14 ("footer") below-section row: fill in date and status like the other rows
$ws.Range("C14").Value = "'25/08"
$ws.Range("E14").Value = "Lựu"
$ws.Range("F14").Value = "Đang làm"
